# Issue-153 Add show pixels checkbox to EPMD
#
# This script reproduces, via the Excel COM object model, the edits that
# were made to the IssuesLog workbook for issue #153:
#   * Several existing OPEN issues (#136, #137, #139, #140, #144) that were
#     related/blocking work get their Status flagged back to "OPEN".
#   * Issue #139 gets a note in the "Sequence" column cross-referencing the
#     new sub-issues (149-152).
#   * Five new issues (149-153) are appended, #153 itself ("Add show pixels
#     checkbox to EPMD") being the one closing out this commit, marked DONE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Re-open a handful of related issues (Status column, G) that this issue
# touches/depends on.
# ---------------------------------------------------------------------------
$ws.Range("G136").Value = "OPEN"
$ws.Range("G137").Value = "OPEN"
$ws.Range("G139").Value = "OPEN"
$ws.Range("G140").Value = "OPEN"
$ws.Range("G144").Value = "OPEN"

# Issue 139 ("Check performance of all Actions in EPMD and their
# draggability") now cross references the four new sub-issues below.
$ws.Range("C139").Value = "149, 150, 151, 152"

# ---------------------------------------------------------------------------
# Append the new issues, rows 149-153.
# ---------------------------------------------------------------------------

# Row 149
$ws.Range("A149").Value = 149
$ws.Range("B149").Value = "Action pixel of and off wide only reapproximates new curves with straight lines"
$ws.Range("C149").Value = -100
$ws.Range("D149").Value = 139
$ws.Range("F149").Value = 43707
$ws.Range("G149").Value = "OPEN"

# Row 150
$ws.Range("A150").Value = 150
$ws.Range("B150").Value = "Toggle needs batch acceleration like pixel, and also approximate with curves"
$ws.Range("D150").Value = 139
$ws.Range("F150").Value = 43707
$ws.Range("G150").Value = "OPEN"

# Row 151
$ws.Range("A151").Value = 151
$ws.Range("B151").Value = "Check clickability of delete pixelChain and the wide options"
$ws.Range("D151").Value = 139
$ws.Range("F151").Value = 43707
$ws.Range("G151").Value = "OPEN"

# Row 152
$ws.Range("A152").Value = 152
$ws.Range("B152").Value = "Change pixelchain thickness should have a draggable option"
$ws.Range("C152").Value = -100
$ws.Range("D152").Value = 139
$ws.Range("F152").Value = 43707
$ws.Range("G152").Value = "OPEN"

# Row 153 - the issue this commit actually closes out.
$ws.Range("A153").Value = 153
$ws.Range("B153").Value = "Add show pixels checkbox to EPMD"
$ws.Range("C153").Value = -100
$ws.Range("F153").Value = 43707
$ws.Range("G153").Value = "DONE"

# ---------------------------------------------------------------------------
# Update the sheet's view/dimension bookkeeping to match the new data extent.
# ---------------------------------------------------------------------------
$ws.Range("A153:B153").Select()
